# chore: adapt column header formatting to respective input file names
#
# The sheet compares an "old" AHB version against a "new" one. Instead of
# the generic "_old"/"_new" suffixes, the header row should be suffixed
# with the concrete format versions being compared: "_FV2404" (old/left
# side) and "_FV2410" (new/right side). The data itself (rows 2-57) is
# untouched. The range is also turned into a proper Excel Table and the
# header row is frozen.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fv2404Headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

$fv2410Headers = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

# Columns A-J: "<field>_FV2404"
for ($i = 0; $i -lt $fv2404Headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value2 = $fv2404Headers[$i]
}

# Column K ("diff") is unchanged.

# Columns L-U: "<field>_FV2410"
for ($i = 0; $i -lt $fv2410Headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 12).Value2 = $fv2410Headers[$i]
}

# Convert the A1:U57 range into a native Excel Table ("Table1") with an
# autofilter on the header row.
$range = $ws.Range("A1:U57")
$tbl = $ws.ListObjects.Add(1, $range, $null, 1)
$tbl.Name = "Table1"

# Freeze the header row (split below row 1, keep the header pane active).
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
